{"js": "// Spanish (es-MX) translation pass for \"Take a Pause\" video script.\n// Replace each English/old-Spanish source string with its translated\n// counterpart. Matching is done with exact (case-sensitive) searches on\n// context.document.body so each run is located unambiguously; duplicate\n// occurrences (e.g. the repeated \"In;\" / \"and out;\" breathing cues) are\n// all replaced via the full result set returned by search().\n\nconst replacements = [\n  [\n    \"Si\u00e9ntate en un lugar c\u00f3modo y cierra los ojos si te parece bien. \",\n    \"Si\u00e9ntate en un lugar c\u00f3modo y, si te parece bien, cierra los ojos. \"\n  ],\n  // \"In;\" is a prefix of the other \"In; \" runs and \"and out;\" is a prefix\n  // of \"and out; \", so a single search+replace on the shorter form covers\n  // every occurrence (3x \"In;\"/\"In; \" and 3x \"and out;\"/\"and out; \")\n  // without the two searches double-matching the same range.\n  [\"In;\", \"Entra;\"],\n  [\"and out;\", \"y sale;\"],\n  [\n    \"Notice how your body feels while you breathe. \",\n    \"F\u00edjate c\u00f3mo se siente tu cuerpo mientras respiras. \"\n  ],\n  [\n    \"Notice where you feel tension in your body and try to let it go.\",\n    \"Nota d\u00f3nde sientes tensi\u00f3n en el cuerpo e intenta soltarla.\"\n  ],\n  [\"Try to let your body relax.\", \"Busca relajar tu cuerpo.\"],\n  [\n    \"When you are ready, open your eyes again. \",\n    \"Cuando sientas que ya has terminado, abre los ojos. \"\n  ],\n  [\n    \"Now, notice if you are feeling any differently than when you started this activity.\",\n    \"Ahora, intenta identificar si te sientes diferente de cuando empezaste esta actividad.\"\n  ],\n  [\n    \"Try slowing down whenever you feel angry, overwhelmed, stressed or worried. \\u23f0 Even a few deep breaths or connecting with the ground beneath you can make a difference. You can also slow down with Your Child!\",\n    \"Intenta ir m\u00e1s lento cada vez que sientas enojo, agobio, estr\u00e9s o preocupaci\u00f3n. \\u23f0 Aunque sea un par de respiraciones profundas o tomarte unos segundos para sentir c\u00f3mo te conectas con el suelo pueden marcar la diferencia. \u00a1Tambi\u00e9n puedes intentar ir m\u00e1s lento con tu ni\u00f1o o ni\u00f1a!\"\n  ],\n  [\n    \"Try slowing down whenever you feel angry, overwhelmed, stressed or worried.  Even a few deep breaths or connecting with the ground beneath you can make a difference. You can also slow down with your girl, boy, or teen!\",\n    \"Intenta ir m\u00e1s lento cada vez que sientas enojo, agobio, estr\u00e9s o preocupaci\u00f3n.  Aunque sea un par de respiraciones profundas o tomarte unos segundos para sentir c\u00f3mo te conectas con el suelo pueden marcar la diferencia. \u00a1Tambi\u00e9n puedes intentar ir m\u00e1s lento con tu ni\u00f1o, ni\u00f1a o adolescente!\"\n  ],\n  [\n    \"Try it with your girl, boy, or teen! \",\n    \"\u00a1Pru\u00e9balo con tu ni\u00f1a, ni\u00f1o o adolescente! \"\n  ]\n];\n\nconst body = context.document.body;\nconst allResults = [];\n\nfor (const [find] of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  allResults.push(found);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, replaceWith] = replacements[i];\n  const found = allResults[i];\n  for (const range of found.items) {\n    range.insertText(replaceWith, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Spanish (es-MX) translation pass for \"Take a Pause\" video script.\n# Uses Find/Replace (\"Replace All\") against the whole document content so\n# every occurrence of each source string is updated, including the\n# repeated \"In;\" / \"and out;\" breathing cues that appear several times.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 2          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$FindText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]2, [ref]$false, [ref]$ReplaceText, [ref]2) | Out-Null\n}\n\nReplace-AllText \"Si\u00e9ntate en un lugar c\u00f3modo y cierra los ojos si te parece bien. \" \"Si\u00e9ntate en un lugar c\u00f3modo y, si te parece bien, cierra los ojos. \"\n\n# \"In;\" is a prefix of the other \"In; \" runs, and \"and out;\" is a prefix\n# of \"and out; \" \u2014 replacing the shorter form covers every occurrence\n# (3x \"In;\"/\"In; \" and 3x \"and out;\"/\"and out; \") in one pass each.\nReplace-AllText \"In;\" \"Entra;\"\nReplace-AllText \"and out;\" \"y sale;\"\n\nReplace-AllText \"Notice how your body feels while you breathe. \" \"F\u00edjate c\u00f3mo se siente tu cuerpo mientras respiras. \"\nReplace-AllText \"Notice where you feel tension in your body and try to let it go.\" \"Nota d\u00f3nde sientes tensi\u00f3n en el cuerpo e intenta soltarla.\"\nReplace-AllText \"Try to let your body relax.\" \"Busca relajar tu cuerpo.\"\nReplace-AllText \"When you are ready, open your eyes again. \" \"Cuando sientas que ya has terminado, abre los ojos. \"\nReplace-AllText \"Now, notice if you are feeling any differently than when you started this activity.\" \"Ahora, intenta identificar si te sientes diferente de cuando empezaste esta actividad.\"\nReplace-AllText \"Try slowing down whenever you feel angry, overwhelmed, stressed or worried. \u23f0 Even a few deep breaths or connecting with the ground beneath you can make a difference. You can also slow down with Your Child!\" \"Intenta ir m\u00e1s lento cada vez que sientas enojo, agobio, estr\u00e9s o preocupaci\u00f3n. \u23f0 Aunque sea un par de respiraciones profundas o tomarte unos segundos para sentir c\u00f3mo te conectas con el suelo pueden marcar la diferencia. \u00a1Tambi\u00e9n puedes intentar ir m\u00e1s lento con tu ni\u00f1o o ni\u00f1a!\"\nReplace-AllText \"Try slowing down whenever you feel angry, overwhelmed, stressed or worried.  Even a few deep breaths or connecting with the ground beneath you can make a difference. You can also slow down with your girl, boy, or teen!\" \"Intenta ir m\u00e1s lento cada vez que sientas enojo, agobio, estr\u00e9s o preocupaci\u00f3n.  Aunque sea un par de respiraciones profundas o tomarte unos segundos para sentir c\u00f3mo te conectas con el suelo pueden marcar la diferencia. \u00a1Tambi\u00e9n puedes intentar ir m\u00e1s lento con tu ni\u00f1o, ni\u00f1a o adolescente!\"\nReplace-AllText \"Try it with your girl, boy, or teen! \" \"\u00a1Pru\u00e9balo con tu ni\u00f1a, ni\u00f1o o adolescente! \"\n"}
